$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (mirrors existing USER NAME / PASSWORD rows)
$ws.Range("A11").Value = "SECR-008"
$ws.Range("B11").Value = "Fosroc@1"
$ws.Range("C11").Value = "Web login SECR"

# Add hyperlink for the new password cell, same pattern as the other rows
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Fosroc@1") | Out-Null

# Match formatting of the other password cells (Hyperlink style, same as column B elsewhere)
$ws.Range("B11").Style = "Hyperlink"

# Update selection to mirror the recorded cursor position after edit
$ws.Range("F11").Select()
